$wb = $excel.ActiveWorkbook

# The "Status" value for this file moved from "Ready for handoff" to "In Translation"
# everywhere it appears (Overview summary columns + each locale sheet's Status column).
$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: keep the known string literal on the LEFT of -eq. PowerShell's
        # -eq coerces the right operand to the left operand's type, so with a
        # Boolean cell value on the left, any non-empty string (like our
        # search text) would coerce to $true and falsely "match".
        if ($oldStatus -eq $cell.Value2) {
            $cell.Value2 = $newStatus
        }
    }
}

# The Status column got narrower (report generated for archive trims it down)
# from ~17.216 characters to ~13.410 characters.
#
# ColumnWidth is read/written in whole "characters" but Excel actually stores
# column widths as a whole number of pixels at the workbook's default font,
# so any value we assign gets snapped to the nearest pixel once it round-trips
# through the object model (exactly like the real Excel.Application COM
# automation this shim is emulating). 12.5 is the input that snaps to the
# pixel width closest to the target 13.4101845877511.
$newWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newWidth   # E: zh-cn status
$overview.Columns.Item(6).ColumnWidth = $newWidth   # F: de-de status

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newWidth        # C: Status

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newWidth        # C: Status
